$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Faculty")

# Remove the now-unused trailing rows (old rows 9-15); final used range is A1:B8
$ws.Range("A9:B15").EntireRow.Delete()

# --- Column A (Name) - Turkish faculty names, written top to bottom -----
$ws.Range("A2").Value = "ECZACILIK FAKÜLTESİ"
$ws.Range("A3").Value = "EĞİTİM BİLİMLERİ FAKÜLTESİ"
$ws.Range("A4").Value = "HUKUK FAKÜLTESİ"
$ws.Range("A5").Value = "İKTİSADİ, İDARİ VE SOSYAL BİLİMLER FAKÜLTESİ"
$ws.Range("A6").Value = "İLETİŞİM FAKÜLTESİ"
$ws.Range("A7").Value = "MÜHENDİSLİK VE DOĞA BİLİMLERİ FAKÜLTESİ"
$ws.Range("A8").Value = "UYGULAMALI BİLİMLER YÜKSEKOKULU"

# --- Column B (Description) - written in the order the source material was
#     pasted in (not top-to-bottom), which matches how the shared-string
#     table ends up ordered in the saved workbook.
$ws.Range("B4").Value = "Faculty dedicated to the study of law, providing education and training to future legal professionals."
$ws.Range("B6").Value = "Faculty that specializes in the study of communication, providing comprehensive education and practical skills in various aspects of media, journalism, public relations, and related fields."
$ws.Range("B5").Value = "Faculty focused on the fields of economics, administration, and social sciences, offering comprehensive education and research opportunities in these disciplines."
$ws.Range("B7").Value = "Faculty dedicated to engineering and natural sciences, offering comprehensive education and research opportunities in fields such as physics, chemistry, mathematics, and various branches of engineering."
$ws.Range("B8").Value = "Higher education institution focused on applied sciences, providing practical-oriented education and training in fields such as technology, computer science, healthcare, and other vocational disciplines."
$ws.Range("B2").Value = "Faculty dedicated to the study of pharmacy, offering comprehensive education and training to future pharmacists, equipping them with knowledge and skills in pharmaceutical sciences, patient care, and medication management."
$ws.Range("B3").Value = "Faculty specializing in the field of education, providing comprehensive education and research opportunities for future educators, administrators, and researchers, with a focus on pedagogy, curriculum development, and educational psychology."

# --- Formatting: every name cell gets top-vertical alignment; most also
#     switch from the default Calibri 11 to Arial 10 (row 3 keeps Calibri).
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10

$ws.Range("A3").VerticalAlignment = -4160

$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10

$ws.Range("A5").VerticalAlignment = -4160
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10

$ws.Range("A6").VerticalAlignment = -4160
$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Size = 10

$ws.Range("A7").VerticalAlignment = -4160
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 10

$ws.Range("A8").VerticalAlignment = -4160
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 10

# --- View / window adjustments -------------------------------------------
$ws.Range("A8").Select()
$excel.ActiveWindow.Zoom = 170
